$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New logbook entry appended as row 7. A7 holds a date-like string
# ("2023-12-11"); Excel would otherwise auto-convert that into a date
# serial number on assignment, so it is entered with a leading
# apostrophe (exactly like typing '2023-12-11 into a General cell)
# to force literal text, then the cell style is reset to Normal so no
# lingering quote-prefix formatting is left on the cell.
$ws.Range("A7").Value = "'2023-12-11"
$ws.Range("A7").Style = "Normal"

$ws.Range("B7").Value = "Final1"
$ws.Range("C7").Value = "Final1"
$ws.Range("D7").Value = "Final1"
$ws.Range("E7").Value = "Final1"
$ws.Range("F7").Value = "Final1"
$ws.Range("G7").Value = "secured"
$ws.Range("H7").Value = "Final1"
$ws.Range("I7").Value = "Final1"
$ws.Range("J7").Value = "Final1-2023-12-11.csv"
